# Apply the "add 2022-Q4 data" edit:
#  1. Insert a brand-new worksheet named "2022-Q4" right after "总计",
#     pushing all the existing quarter sheets one slot to the right
#     (their own names/content are untouched).
#  2. Populate the new sheet with the 2022-Q4 per-fund holdings table.
#  3. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q4
#     and shift the previous rows down by one (their values are unchanged).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right before the current "2022-Q3"
#    sheet (i.e. right after "总计").
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

# Match the page-margin metadata used by the sibling quarter sheets.
$ps = $q4.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Fill in the 2022-Q4 holdings table.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# index, code, name, size, position%, ratio%, marketvalue, rank
$rows = @(
    @(0,  "001071", "华安媒体互联网混合A",          "51.25", "91.12", "3.01", "1.5426", 8),
    @(1,  "017766", "华夏兴和混合C",                "35.60", "91.76", "3.37", "1.1997", 10),
    @(2,  "210009", "金鹰核心资源混合",              "3.14",  "91.78", "6.79", "0.2132", 1),
    @(3,  "001167", "金鹰科技创新股票",              "3.17",  "91.02", "6.49", "0.2057", 1),
    @(4,  "162102", "金鹰中小盘精选混合",            "3.48",  "76.23", "5.43", "0.1890", 3),
    @(5,  "013620", "华安媒体互联网混合C",          "2.39",  "91.12", "3.01", "0.0719", 8),
    @(6,  "210002", "金鹰红利价值混合A",            "0.96",  "77.22", "6.54", "0.0628", 1),
    @(7,  "290012", "泰信行业精选灵活配置混合A",    "1.17",  "92.06", "5.08", "0.0594", 5),
    @(8,  "159855", "银华中证影视主题ETF",          "1.01",  "97.80", "4.66", "0.0471", 5),
    @(9,  "016563", "金鹰红利价值混合C",            "0.52",  "77.22", "6.54", "0.0340", 1),
    @(10, "516620", "国泰中证影视主题ETF",          "0.71",  "98.01", "4.28", "0.0304", 7),
    @(11, "002583", "泰信行业精选灵活配置混合C",    "0.23",  "92.06", "5.08", "0.0117", 5),
    @(12, "004677", "博时战略新兴产业混合",          "0.34",  "86.55", "3.27", "0.0111", 10),
    @(13, "001613", "长城久祥灵活配置混合A",        "0.24",  "94.17", "4.51", "0.0108", 4),
    @(14, "004988", "人保双利优选混合A",            "0.57",  "29.15", "0.51", "0.0029", 6),
    @(15, "004989", "人保双利优选混合C",            "0.01",  "29.15", "0.51", "0.0001", 6),
    @(16, "017462", "长城久祥灵活配置混合C",        "0.00",  "94.17", "4.51", $null,    4),
    @(17, "519918", "华夏兴和混合A",                "0.00",  "91.76", "3.37", $null,    10)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    if ($row[6] -eq $null) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    }
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: shift rows 2-8 down to 3-9 (copy,
#    which preserves the original formatting), then overwrite row 2
#    with the new 2022-Q4 totals.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Range("A2:D8").Copy($summary.Range("A3"))

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 18
$summary.Cells.Item(2, 4).Value = 3.69

# Renumber the A-column index (0-based row counter) for the shifted rows.
for ($i = 3; $i -le 9; $i++) {
    $summary.Cells.Item($i, 1).Value = $i - 2
}

# Keep "总计" as the active sheet/tab, matching the original workbook.
$summary.Activate()
